$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TIEMPO_CAJA_NORMAL")

$ws.Cells.Item(2, 1).Value = "00:02:51"
$ws.Cells.Item(3, 1).Value = "00:03:09"
$ws.Cells.Item(4, 1).Value = "00:01:03"
$ws.Cells.Item(5, 1).Value = "00:05:45"
$ws.Cells.Item(6, 1).Value = "00:00:22"
$ws.Cells.Item(7, 1).Value = "00:05:42"
$ws.Cells.Item(8, 1).Value = "00:01:15"
$ws.Cells.Item(9, 1).Value = "00:02:52"
$ws.Cells.Item(10, 1).Value = "00:02:47"
$ws.Cells.Item(11, 1).Value = "00:07:56"
$ws.Cells.Item(12, 1).Value = "00:07:14"
$ws.Cells.Item(13, 1).Value = "00:00:04"
$ws.Cells.Item(14, 1).Value = "00:00:01"
$ws.Cells.Item(15, 1).Value = "00:01:34"
$ws.Cells.Item(16, 1).Value = "00:00:04"
$ws.Cells.Item(17, 1).Value = "00:02:08"
$ws.Cells.Item(18, 1).Value = "00:01:31"
$ws.Cells.Item(19, 1).Value = "00:02:11"
$ws.Cells.Item(20, 1).Value = "00:01:09"
$ws.Cells.Item(21, 1).Value = "00:00:38"
$ws.Cells.Item(22, 1).Value = "00:05:58"
$ws.Cells.Item(23, 1).Value = "00:02:36"
$ws.Cells.Item(24, 1).Value = "00:00:00"
$ws.Cells.Item(25, 1).Value = "00:00:12"
$ws.Cells.Item(26, 1).Value = "00:01:59"
$ws.Cells.Item(27, 1).Value = "00:00:40"
$ws.Cells.Item(28, 1).Value = "00:01:02"
$ws.Cells.Item(29, 1).Value = "00:03:02"
$ws.Cells.Item(30, 1).Value = "00:00:52"
$ws.Cells.Item(31, 1).Value = "00:05:12"
$ws.Cells.Item(32, 1).Value = "00:02:32"
$ws.Cells.Item(33, 1).Value = "00:01:36"
$ws.Cells.Item(34, 1).Value = "00:00:23"
$ws.Cells.Item(35, 1).Value = "00:01:03"
$ws.Cells.Item(36, 1).Value = "00:02:52"
$ws.Cells.Item(37, 1).Value = "00:09:22"
$ws.Cells.Item(38, 1).Value = "00:00:09"
$ws.Cells.Item(39, 1).Value = "00:00:17"
$ws.Cells.Item(40, 1).Value = "00:02:41"
$ws.Cells.Item(41, 1).Value = "00:03:35"
$ws.Cells.Item(42, 1).Value = "00:00:11"
$ws.Cells.Item(43, 1).Value = "00:03:32"
$ws.Cells.Item(44, 1).Value = "00:03:05"
$ws.Cells.Item(45, 1).Value = "00:02:19"
$ws.Cells.Item(46, 1).Value = "00:00:49"
$ws.Cells.Item(47, 1).Value = "00:01:10"
$ws.Cells.Item(48, 1).Value = "00:00:54"
$ws.Cells.Item(49, 1).Value = "00:03:40"
$ws.Cells.Item(50, 1).Value = "00:00:47"
$ws.Cells.Item(51, 1).Value = "00:00:46"
$ws.Cells.Item(52, 1).Value = "00:05:28"
$ws.Cells.Item(53, 1).Value = "00:05:03"
$ws.Cells.Item(54, 1).Value = "00:06:55"
$ws.Cells.Item(55, 1).Value = "00:03:17"
$ws.Cells.Item(56, 1).Value = "00:00:10"
$ws.Cells.Item(57, 1).Value = "00:00:11"
$ws.Cells.Item(58, 1).Value = "00:01:37"
$ws.Cells.Item(59, 1).Value = "00:03:52"
$ws.Cells.Item(60, 1).Value = "00:04:18"
$ws.Cells.Item(61, 1).Value = "00:02:01"
$ws.Cells.Item(62, 1).Value = "00:06:14"
$ws.Cells.Item(63, 1).Value = "00:00:19"
$ws.Cells.Item(64, 1).Value = "00:00:30"
$ws.Cells.Item(65, 1).Value = "00:00:48"
$ws.Cells.Item(66, 1).Value = "00:00:35"
$ws.Cells.Item(67, 1).Value = "00:00:44"
$ws.Cells.Item(68, 1).Value = "00:00:15"
$ws.Cells.Item(69, 1).Value = "00:00:56"
$ws.Cells.Item(70, 1).Value = "00:00:03"
$ws.Cells.Item(71, 1).Value = "00:01:32"
$ws.Cells.Item(72, 1).Value = "00:05:49"
$ws.Cells.Item(73, 1).Value = "00:01:10"
$ws.Cells.Item(74, 1).Value = "00:01:57"
$ws.Cells.Item(75, 1).Value = "00:01:53"
$ws.Cells.Item(76, 1).Value = "00:05:25"
$ws.Cells.Item(77, 1).Value = "00:02:39"
$ws.Cells.Item(78, 1).Value = "00:03:29"
$ws.Cells.Item(79, 1).Value = "00:00:49"
$ws.Cells.Item(80, 1).Value = "00:04:45"
$ws.Cells.Item(81, 1).Value = "00:00:24"
$ws.Cells.Item(82, 1).Value = "00:01:02"
$ws.Cells.Item(83, 1).Value = "00:02:40"
$ws.Cells.Item(84, 1).Value = "00:00:48"
$ws.Cells.Item(85, 1).Value = "00:04:20"
$ws.Cells.Item(86, 1).Value = "00:03:56"
$ws.Cells.Item(87, 1).Value = "00:00:21"
$ws.Cells.Item(88, 1).Value = "00:00:31"
$ws.Cells.Item(89, 1).Value = "00:00:06"
$ws.Cells.Item(90, 1).Value = "00:03:54"
$ws.Cells.Item(91, 1).Value = "00:00:49"
$ws.Cells.Item(92, 1).Value = "00:00:59"
$ws.Cells.Item(93, 1).Value = "00:00:37"
$ws.Cells.Item(94, 1).Value = "00:00:15"
$ws.Cells.Item(95, 1).Value = "00:01:13"
$ws.Cells.Item(96, 1).Value = "00:00:32"
$ws.Cells.Item(97, 1).Value = "00:01:39"
$ws.Cells.Item(98, 1).Value = "00:00:27"
$ws.Cells.Item(99, 1).Value = "00:01:25"
$ws.Cells.Item(100, 1).Value = "00:05:55"
$ws.Cells.Item(101, 1).Value = "00:00:55"
